# Update odds values on Sheet1 (Jogos_da_Semana_FlashScore_2024-11-26)
# for rows 10, 14, 17, 18, 20, 21 to match the latest FlashScore data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G10").Value = 1.67
$ws.Range("I10").Value = 4.5
$ws.Range("AD10").Value = 8
$ws.Range("AU10").Value = 7
$ws.Range("G14").Value = 4.33
$ws.Range("I14").Value = 1.7
$ws.Range("J14").Value = 4.75
$ws.Range("L14").Value = 2.3
$ws.Range("AB14").Value = 34
$ws.Range("AX14").Value = 9
$ws.Range("AZ14").Value = 29
$ws.Range("BC14").Value = 501
$ws.Range("G17").Value = 2
$ws.Range("I17").Value = 3.75
$ws.Range("Q17").Value = 2.1
$ws.Range("R17").Value = 1.7
$ws.Range("W17").Value = 6.5
$ws.Range("X17").Value = 9
$ws.Range("Z17").Value = 17
$ws.Range("AA17").Value = 17
$ws.Range("AH17").Value = 10
$ws.Range("AI17").Value = 19
$ws.Range("AV17").Value = 67
$ws.Range("AY17").Value = 34
$ws.Range("G18").Value = 1.65
$ws.Range("H18").Value = 4
$ws.Range("I18").Value = 4.5
$ws.Range("J18").Value = 2.25
$ws.Range("W18").Value = 8
$ws.Range("X18").Value = 8.5
$ws.Range("Z18").Value = 13
$ws.Range("AD18").Value = 7.5
$ws.Range("AL18").Value = 34
$ws.Range("AN18").Value = 3.75
$ws.Range("AO18").Value = 8.5
$ws.Range("AW18").Value = 6.5
$ws.Range("G20").Value = 5.25
$ws.Range("H20").Value = 4.2
$ws.Range("I20").Value = 1.55
$ws.Range("K20").Value = 2.4
$ws.Range("N20").Value = 15
$ws.Range("O20").Value = 1.18
$ws.Range("P20").Value = 4.5
$ws.Range("Q20").Value = 1.6
$ws.Range("R20").Value = 2.3
$ws.Range("S20").Value = 1.3
$ws.Range("T20").Value = 3.4
$ws.Range("W20").Value = 17
$ws.Range("X20").Value = 29
$ws.Range("AC20").Value = 15
$ws.Range("AH20").Value = 9
$ws.Range("AL20").Value = 12
$ws.Range("AN20").Value = 7
$ws.Range("AT20").Value = 3.4
$ws.Range("AV20").Value = 51
$ws.Range("AY20").Value = 17
$ws.Range("G21").Value = 1.42
$ws.Range("I21").Value = 7
$ws.Range("J21").Value = 1.91
$ws.Range("N21").Value = 15
$ws.Range("W21").Value = 8.5
$ws.Range("Z21").Value = 10
$ws.Range("AC21").Value = 15
$ws.Range("AE21").Value = 17
$ws.Range("AF21").Value = 51
$ws.Range("AG21").Value = 201
$ws.Range("AI21").Value = 41
$ws.Range("AJ21").Value = 21
$ws.Range("AK21").Value = 81
$ws.Range("AP21").Value = 17
$ws.Range("AU21").Value = 8.5
$ws.Range("AX21").Value = 34
$ws.Range("BA21").Value = 126
